$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the original "bkash number" numeric cell style (s=3, quotePrefix)
# from F2 onto a scratch cell far outside the used range, so we can restore
# it on F2:F5 after their values are overwritten below.
$ws.Range("F2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)

# --- Row 2: update existing retailer entry (Mugdho Corporation stays) ---
$ws.Range("B2").Value = "RET-08787"
$ws.Range("C2").Value = "Shabbir Telecom"
$ws.Range("D2").Value = "Md Santu Ali"
$ws.Range("F2").Value = 1711015278

# Fix B2's cell style: it previously used a bold header-ish style (s=4);
# the new layout uses the plain bordered style shared by the rest of the
# data rows (same style as A2/C2/D2/E2).
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# --- Row 3: new retailer entry ---
$ws.Range("A3").Value = "Mugdho Corporation"
$ws.Range("B3").Value = "RET-07875"
$ws.Range("C3").Value = "Khondokar Elecrtronics"
$ws.Range("D3").Value = "Md Monirul Islam"
$ws.Range("E3").Value = "Bkash "
$ws.Range("F3").Value = 1611438268

# --- Row 4: new retailer entry (note: D filled before C to match source
# workbook's shared-string ordering) ---
$ws.Range("A4").Value = "Mugdho Corporation"
$ws.Range("B4").Value = "RET-35442"
$ws.Range("D4").Value = "Md Mojnu Pramanic"
$ws.Range("C4").Value = "Moom Telecom"
$ws.Range("E4").Value = "Bkash "
$ws.Range("F4").Value = 1712469447
# C4 keeps the default (no border) style, unlike its row neighbours.
$ws.Range("C4").ClearFormats()

# --- Row 5: new retailer entry ---
$ws.Range("A5").Value = "Mugdho Corporation"
$ws.Range("B5").Value = "RET-14872"
$ws.Range("C5").Value = "Bismillah Mobile Center"
$ws.Range("D5").Value = "Md. Sohel Mandal"
$ws.Range("E5").Value = "Bkash "
$ws.Range("F5").Value = 1738182165

# Restore the stashed numeric "bkash number" style onto F2:F5 (setting
# .Value above reset their style to the plain bordered style).
$ws.Range("Z100").Copy()
$ws.Range("F2:F5").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# --- Column widths (widened for the longer retailer/shop names) ---
$ws.Columns("C").ColumnWidth = 21.833333333333332
$ws.Columns("D").ColumnWidth = 18

# --- Selection moved from K19 to G19 ---
$ws.Range("G19").Select() | Out-Null
